$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H53").Value = 691.5
$ws.Range("I53").Value = 850
$ws.Range("J53").Value = 659.8
$ws.Range("K53").Value = 850
$ws.Range("L53").Value = 659.8
$ws.Range("M53").Value = -213
$ws.Range("N53").Value = -1933.8
$ws.Range("H64").Value = 21973
$ws.Range("I64").Value = 28959.75
$ws.Range("J64").Value = 7999.5
$ws.Range("K64").Value = 28959.75
$ws.Range("L64").Value = 7999.5
$ws.Range("M64").Value = -28711.75
$ws.Range("N64").Value = -8495.5
$ws.Range("H67").Value = 21973
$ws.Range("I67").Value = 28959.75
$ws.Range("J67").Value = 7999.5
$ws.Range("K67").Value = 28959.75
$ws.Range("L67").Value = 7999.5
$ws.Range("M67").Value = -28101.75
$ws.Range("N67").Value = -9715.5
$ws.Range("H92").Value = 2227.72
$ws.Range("I92").Value = 931.5625
$ws.Range("J92").Value = 4532
$ws.Range("K92").Value = 931.5625
$ws.Range("L92").Value = 4532
$ws.Range("M92").Value = 316.4375
$ws.Range("N92").Value = -7028
$ws.Range("H125").Value = 5507651.5
$ws.Range("I125").Value = 15907453
$ws.Range("K125").Value = 143167077
$ws.Range("M125").Value = -143164617
$ws.Range("H137").Value = 10401.389
$ws.Range("I137").Value = 1778.8572
$ws.Range("J137").Value = 15888.454
$ws.Range("K137").Value = 5336.571599999999
$ws.Range("L137").Value = 47665.362
$ws.Range("M137").Value = -2786.571599999999
$ws.Range("N137").Value = -52765.362
$ws.Range("H138").Value = 1419834.2
$ws.Range("J138").Value = 2054026.9
$ws.Range("L138").Value = 6162080.699999999
$ws.Range("N138").Value = -6172360.699999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18369.033
$ws.Range("I32").Value = 20482.703
$ws.Range("K32").Value = 20482.703
$ws.Range("M32").Value = -20195.703
$ws.Range("H61").Value = 6483.2285
$ws.Range("I61").Value = 3720.9656
$ws.Range("K61").Value = 3720.9656
$ws.Range("M61").Value = -3508.9656
$ws.Range("H74").Value = 4213.523
$ws.Range("I74").Value = 1766
$ws.Range("J74").Value = 7150.55
$ws.Range("K74").Value = 1766
$ws.Range("L74").Value = 7150.55
$ws.Range("M74").Value = -892
$ws.Range("N74").Value = -8898.549999999999
$ws.Range("H77").Value = 4213.523
$ws.Range("I77").Value = 1766
$ws.Range("J77").Value = 7150.55
$ws.Range("K77").Value = 8830
$ws.Range("L77").Value = 35752.75
$ws.Range("M77").Value = -4462
$ws.Range("N77").Value = -44488.75
$ws.Range("H110").Value = 24347.352
$ws.Range("I110").Value = 30105.414
$ws.Range("K110").Value = 30105.414
$ws.Range("M110").Value = -28060.414
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").ClearContents()
$ws.Range("N112").Value = 0
$ws.Range("H132").Value = 1479.5968
$ws.Range("I132").Value = 1241.2037
$ws.Range("J132").Value = 3088.75
$ws.Range("K132").Value = 3723.6111
$ws.Range("L132").Value = 9266.25
$ws.Range("M132").Value = -1193.6111
$ws.Range("N132").Value = -14326.25
$ws.Range("H136").Value = 6483.2285
$ws.Range("I136").Value = 3720.9656
$ws.Range("K136").Value = 11162.8968
$ws.Range("M136").Value = -8612.8968
$ws.Range("H139").Value = 130698.336
$ws.Range("J139").Value = 130698.336
$ws.Range("L139").Value = 130698.336
$ws.Range("N139").Value = -140978.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1088.2307
$ws.Range("I107").Value = 887.25
$ws.Range("K107").Value = 887.25
$ws.Range("M107").Value = 1032.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 10600611
$ws.Range("I19").Value = 11564249
$ws.Range("J19").Value = 596
$ws.Range("K19").Value = 11564249
$ws.Range("L19").Value = 596
$ws.Range("M19").Value = -11564079
$ws.Range("N19").Value = -936
$ws.Range("H24").Value = 10600611
$ws.Range("I24").Value = 11564249
$ws.Range("J24").Value = 596
$ws.Range("K24").Value = 11564249
$ws.Range("L24").Value = 596
$ws.Range("M24").Value = -11564079
$ws.Range("N24").Value = -936
$ws.Range("H29").Value = 25500
$ws.Range("J29").Value = 26000
$ws.Range("L29").Value = 26000
$ws.Range("N29").Value = -26586
$ws.Range("H31").Value = 1925570.5
$ws.Range("I31").Value = 4764345
$ws.Range("J31").Value = 2529.742
$ws.Range("K31").Value = 4764345
$ws.Range("L31").Value = 2529.742
$ws.Range("M31").Value = -4764050
$ws.Range("N31").Value = -3119.742
$ws.Range("H34").Value = 1925570.5
$ws.Range("I34").Value = 4764345
$ws.Range("J34").Value = 2529.742
$ws.Range("K34").Value = 4764345
$ws.Range("L34").Value = 2529.742
$ws.Range("M34").Value = -4764143
$ws.Range("N34").Value = -2933.742
$ws.Range("H94").Value = 1180.1818
$ws.Range("I94").Value = 1406.75
$ws.Range("J94").Value = 1050.7142
$ws.Range("K94").Value = 1406.75
$ws.Range("L94").Value = 1050.7142
$ws.Range("M94").Value = -955.75
$ws.Range("N94").Value = -1952.7142
$ws.Range("H105").Value = 1599
$ws.Range("I105").Value = 925.8
$ws.Range("K105").Value = 925.8
$ws.Range("M105").Value = 821.2
$ws.Range("H132").Value = 2781.2942
$ws.Range("I132").Value = 2023.25
$ws.Range("K132").Value = 6069.75
$ws.Range("M132").Value = -3539.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 959.4
$ws.Range("J5").Value = 898.5
$ws.Range("L5").Value = 2695.5
$ws.Range("N5").Value = -2919.5
$ws.Range("H17").Value = 1092.6666
$ws.Range("J17").Value = 489
$ws.Range("L17").Value = 1467
$ws.Range("N17").Value = -1805
$ws.Range("H34").Value = 5472077
$ws.Range("I34").Value = 7295958
$ws.Range("J34").Value = 433.33334
$ws.Range("K34").Value = 21887874
$ws.Range("L34").Value = 1300.00002
$ws.Range("M34").Value = -21887790
$ws.Range("N34").Value = -1468.00002
$ws.Range("H39").Value = 10853.333
$ws.Range("J39").Value = 4000
$ws.Range("L39").Value = 12000
$ws.Range("N39").Value = -12588
$ws.Range("H55").Value = 1998.8
$ws.Range("J55").Value = 2000
$ws.Range("L55").Value = 6000
$ws.Range("N55").Value = -6354
$ws.Range("H121").Value = 753971.7
$ws.Range("I121").Value = 143381.86
$ws.Range("K121").Value = 430145.58
$ws.Range("M121").Value = -428835.58
$ws.Range("H135").Value = 959.4
$ws.Range("J135").Value = 898.5
$ws.Range("L135").Value = 8086.5
$ws.Range("N135").Value = -13156.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3156.8462
$ws.Range("I80").Value = 2628.7646
$ws.Range("J80").Value = 4154.3335
$ws.Range("K80").Value = 2628.7646
$ws.Range("L80").Value = 4154.3335
$ws.Range("M80").Value = -1630.7646
$ws.Range("N80").Value = -6150.3335
$ws.Range("H83").Value = 3156.8462
$ws.Range("I83").Value = 2628.7646
$ws.Range("J83").Value = 4154.3335
$ws.Range("K83").Value = 13143.823
$ws.Range("L83").Value = 20771.6675
$ws.Range("M83").Value = -8151.823
$ws.Range("N83").Value = -30755.6675
$ws.Range("H97").Value = 601.8857400000001
$ws.Range("I97").Value = 454.2381
$ws.Range("K97").Value = 454.2381
$ws.Range("M97").Value = 41.76190000000003
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").ClearContents()
$ws.Range("N111").Value = 0
$ws.Range("H126").Value = 3608.2144
$ws.Range("I126").Value = 3270.2222
$ws.Range("J126").Value = 4216.6
$ws.Range("K126").Value = 9810.6666
$ws.Range("L126").Value = 12649.8
$ws.Range("M126").Value = -7340.6666
$ws.Range("N126").Value = -17589.8
$ws.Range("H132").Value = 4150.857
$ws.Range("I132").Value = 4867.375
$ws.Range("J132").Value = 3709.923
$ws.Range("K132").Value = 14602.125
$ws.Range("L132").Value = 11129.769
$ws.Range("M132").Value = -12072.125
$ws.Range("N132").Value = -16189.769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3885
$ws.Range("I46").Value = 1989
$ws.Range("K46").Value = 1989
$ws.Range("M46").Value = -1801
$ws.Range("H55").Value = 1150.7297
$ws.Range("I55").Value = 895.64703
$ws.Range("J55").Value = 1367.55
$ws.Range("K55").Value = 895.64703
$ws.Range("L55").Value = 1367.55
$ws.Range("M55").Value = -722.64703
$ws.Range("N55").Value = -1713.55
$ws.Range("H108").Value = 70000
$ws.Range("J108").Value = 70000
$ws.Range("L108").Value = 70000
$ws.Range("N108").Value = -77680
$ws.Range("H122").Value = 1787.875
$ws.Range("I122").Value = 1686.1428
$ws.Range("K122").Value = 5058.428400000001
$ws.Range("M122").Value = -2608.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6092.52
$ws.Range("I81").Value = 6431.278
$ws.Range("J81").Value = 5221.4287
$ws.Range("K81").Value = 12862.556
$ws.Range("L81").Value = 10442.8574
$ws.Range("M81").Value = -11801.556
$ws.Range("N81").Value = -12564.8574
$ws.Range("H84").Value = 6092.52
$ws.Range("I84").Value = 6431.278
$ws.Range("J84").Value = 5221.4287
$ws.Range("K84").Value = 64312.78
$ws.Range("L84").Value = 52214.287
$ws.Range("M84").Value = -59008.78
$ws.Range("N84").Value = -62822.287
$ws.Range("H113").Value = 779.5161000000001
$ws.Range("I113").Value = 826.2308
$ws.Range("K113").Value = 2478.6924
$ws.Range("M113").Value = -308.6923999999999
$ws.Range("H126").Value = 2874.8572
$ws.Range("I126").Value = 1171.6364
$ws.Range("K126").Value = 3514.9092
$ws.Range("M126").Value = -1044.9092
$ws.Range("H132").Value = 20569.262
$ws.Range("I132").Value = 24519.406
$ws.Range("J132").Value = 4329.778
$ws.Range("K132").Value = 73558.21799999999
$ws.Range("L132").Value = 12989.334
$ws.Range("M132").Value = -71028.21799999999
$ws.Range("N132").Value = -18049.334
